$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = -0.9132452695465898
$ws.Range("K3").Value = -5.194930382006427
$ws.Range("K4").Value = -0.5694394858978932
$ws.Range("K5").Value = -0.0649767940021416
$ws.Range("K6").Value = -0.2067118886112105
